$wb = $excel.ActiveWorkbook

# Rename first sheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Squier"

# Select a cell on sheet1 to match new selection
$ws1.Range("D11").Select()

# Second sheet - add new data row
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = "Exlda"
$ws2.Range("B1").Value = 1000
$ws2.Range("C1").Value = 2001
$ws2.Range("E1").Value = "Olha"
$ws2.Range("D1").Value = "Superstrat"
$ws2.Range("F1").Value = "Acoustic"

$ws2.Range("D1").Select()

$ws1.Select()
$ws1.Range("D11").Select()
